$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new values for Scenario4 row (row 5)
$ws.Range("C5").Value2 = "[N/A; ReLU; Sigmoid]"
$ws.Range("D5").Value2 = "^"

# Update existing cell: Scenario1 activation functions description
$ws.Range("C2").Value2 = "[N/A; Sigmoid; Sigmoid]"

# Add new row 6 for Scenario5
$ws.Range("A6").Value2 = "NeuralNetworkScenario5.mat"
$ws.Range("B6").Value2 = "^"
$ws.Range("C6").Value2 = "^"
$ws.Range("D6").Value2 = "Partially trained network"
$ws.Range("E6").Value2 = "only about a 50% accuracy"

# Update selection to match target state
$ws.Range("G14").Select()
